# Fix Login test suite
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "login" (sheet1): insert a new "${xpath}" column before the existing
# output/message column, and populate it with the XPath locators used by the
# new assertions.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("login")

# Insert a new column at D; everything currently in D:G shifts right to E:H.
$ws1.Columns.Item(4).Insert(-4161)
$ws1.Columns.Item(4).ColumnWidth = 53.5546875

$ws1.Range("D1").Value = '${xpath}'

$ws1.Range("D2").Value = '//span[@id="UserNameRequired"]'
$ws1.Range("D3").Value = '//span[@id="PasswordRequired"]'
$ws1.Range("D4").Value = '//div[@class="text-danger"]'
$ws1.Range("D5").Value = '//div[@class="text-danger"]'
$ws1.Range("D6").Value = '//span[@id="Header_LoginViewHeader_LoginName1"]'

$ws1.Range("D2:D6").VerticalAlignment = -4108

# The old message column (now column E) becomes a wildcard/xpath-matched
# column: clear its old text-wrap formatting and vertically center it.
$ws1.Range("E2:E3").Value = '*'

$ws1.Range("E4:E5").NumberFormat = "General"
$ws1.Range("E4:E5").WrapText = $false

$ws1.Range("E2:E6").VerticalAlignment = -4108

# New trailing blank cell created by the shift (was G2, now H2).
$ws1.Range("H2").Value = ""

$ws1.Range("B3").Select()

# ---------------------------------------------------------------------------
# Sheet "paging" (sheet2): no data changed, just the remembered selection.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("paging")
$ws2.Range("A30").Select()
